$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 277715
$ws.Range("D2").Value = 355121162
$ws.Range("C8").Value = 769
$ws.Range("D8").Value = 1135394
$ws.Range("C10").Value = 105988
$ws.Range("D10").Value = 155584022
$ws.Range("C12").Value = 52214
$ws.Range("D12").Value = 75473447
$ws.Range("C16").Value = 3429
$ws.Range("D16").Value = 4874603
$ws.Range("C20").Value = 4908
$ws.Range("D20").Value = 6874088
$ws.Range("C22").Value = 68858
$ws.Range("D22").Value = 86478690
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 15000
$ws.Range("C27").Value = 257
$ws.Range("D27").Value = 370718
$ws.Range("C28").Value = 29640
$ws.Range("D28").Value = 43449773
$ws.Range("C30").Value = 10179
$ws.Range("D30").Value = 14687443
$ws.Range("C33").Value = 1351
$ws.Range("D33").Value = 1900715
$ws.Range("C35").Value = 1404
$ws.Range("D35").Value = 1979703
$ws.Range("C36").Value = 87144
$ws.Range("D36").Value = 110495459
$ws.Range("C37").Value = 57
$ws.Range("D37").Value = 65593
$ws.Range("C38").Value = 74
$ws.Range("D38").Value = 104153
$ws.Range("C42").Value = 804
$ws.Range("D42").Value = 1183958
$ws.Range("C43").Value = 8
$ws.Range("D43").Value = 12000
$ws.Range("C44").Value = 40775
$ws.Range("D44").Value = 59878093
$ws.Range("C46").Value = 8155
$ws.Range("D46").Value = 11713027
$ws.Range("C48").Value = 1186
$ws.Range("D48").Value = 1647017
$ws.Range("C51").Value = 1729
$ws.Range("D51").Value = 2411827
$ws.Range("C52").Value = 60536
$ws.Range("D52").Value = 76225032
$ws.Range("C56").Value = 340
$ws.Range("D56").Value = 500909
$ws.Range("C58").Value = 25644
$ws.Range("D58").Value = 37651295
$ws.Range("C61").Value = 9837
$ws.Range("D61").Value = 14243937
$ws.Range("C62").Value = 1198
$ws.Range("D62").Value = 1674607
$ws.Range("C65").Value = 1115
$ws.Range("D65").Value = 1566576
$ws.Range("C67").Value = 17776
$ws.Range("D67").Value = 23270436
$ws.Range("C71").Value = 6465
$ws.Range("D71").Value = 9459992
$ws.Range("C72").Value = 4407
$ws.Range("D72").Value = 6405710
$ws.Range("C73").Value = 405
$ws.Range("D73").Value = 568758
$ws.Range("C74").Value = 219
$ws.Range("D74").Value = 308056
$ws.Range("C75").Value = 123920
$ws.Range("D75").Value = 155152951
$ws.Range("C78").Value = 18
$ws.Range("D78").Value = 27000
$ws.Range("C79").Value = 366
$ws.Range("D79").Value = 534782
$ws.Range("C81").Value = 57653
$ws.Range("D81").Value = 84645609
$ws.Range("C82").Value = 69
$ws.Range("D82").Value = 102082
$ws.Range("C84").Value = 26414
$ws.Range("D84").Value = 38250161
$ws.Range("C86").Value = 2305
$ws.Range("D86").Value = 3324328
$ws.Range("C87").Value = 2100
$ws.Range("D87").Value = 2956207
$ws.Range("C88").Value = 23784
$ws.Range("D88").Value = 32297817
$ws.Range("C92").Value = 6356
$ws.Range("D92").Value = 9368176
$ws.Range("C94").Value = 5523
$ws.Range("D94").Value = 8006034
$ws.Range("C96").Value = 386
$ws.Range("D96").Value = 545809
$ws.Range("C97").Value = 321
$ws.Range("D97").Value = 464161
$ws.Range("C98").Value = 5394
$ws.Range("D98").Value = 7473286
$ws.Range("C100").Value = 1368
$ws.Range("D100").Value = 2014062
$ws.Range("C102").Value = 1871
$ws.Range("D102").Value = 2730601
$ws.Range("C104").Value = 65
$ws.Range("D104").Value = 93020
$ws.Range("C105").Value = 95
$ws.Range("D105").Value = 133789
$ws.Range("C106").Value = 125128
$ws.Range("D106").Value = 155153237
$ws.Range("C110").Value = 864
$ws.Range("D110").Value = 1270632
$ws.Range("C112").Value = 48250
$ws.Range("D112").Value = 70826657
$ws.Range("C113").Value = 79
$ws.Range("D113").Value = 117659
$ws.Range("C114").Value = 23819
$ws.Range("D114").Value = 34533935
$ws.Range("C115").Value = 1102
$ws.Range("D115").Value = 1502832
$ws.Range("C118").Value = 1709
$ws.Range("D118").Value = 2404872
$ws.Range("C120").Value = 37421
$ws.Range("D120").Value = 50207259
$ws.Range("C121").Value = 23
$ws.Range("D121").Value = 32653
$ws.Range("C126").Value = 12429
$ws.Range("D126").Value = 18273127
$ws.Range("C127").Value = 3247
$ws.Range("D127").Value = 4688941
$ws.Range("C130").Value = 326
$ws.Range("D130").Value = 472036
$ws.Range("C132").Value = 295
$ws.Range("D132").Value = 413190
$ws.Range("C133").Value = 14150
$ws.Range("D133").Value = 18777683
$ws.Range("C137").Value = 6071
$ws.Range("D137").Value = 8848752
$ws.Range("C139").Value = 4020
$ws.Range("D139").Value = 5799615
$ws.Range("C141").Value = 212
$ws.Range("D141").Value = 289433
$ws.Range("C142").Value = 173
$ws.Range("D142").Value = 250664
$ws.Range("C144").Value = 8872
$ws.Range("D144").Value = 12805026
$ws.Range("C145").Value = 1309
$ws.Range("D145").Value = 1945078
$ws.Range("C146").Value = 180
$ws.Range("D146").Value = 265671
$ws.Range("C147").Value = 26
$ws.Range("D147").Value = 38690
$ws.Range("C149").Value = 77620
$ws.Range("D149").Value = 97513434
$ws.Range("C154").Value = 589
$ws.Range("D154").Value = 870037
$ws.Range("C156").Value = 31064
$ws.Range("D156").Value = 45606600
$ws.Range("C158").Value = 11501
$ws.Range("D158").Value = 16637592
$ws.Range("C160").Value = 1055
$ws.Range("D160").Value = 1482396
$ws.Range("C162").Value = 1253
$ws.Range("D162").Value = 1774226
$ws.Range("C164").Value = 209991
$ws.Range("D164").Value = 261928662
$ws.Range("C166").Value = 151
$ws.Range("D166").Value = 218453
$ws.Range("C170").Value = 788
$ws.Range("D170").Value = 1158790
$ws.Range("C172").Value = 79169
$ws.Range("D172").Value = 116157136
$ws.Range("C175").Value = 29377
$ws.Range("D175").Value = 42316015
$ws.Range("C178").Value = 4418
$ws.Range("D178").Value = 6310386
$ws.Range("C181").Value = 3722
$ws.Range("D181").Value = 5178659
$ws.Range("C184").Value = 229021
$ws.Range("D184").Value = 284180113
$ws.Range("C185").Value = 138
$ws.Range("D185").Value = 151352
$ws.Range("C186").Value = 219
$ws.Range("D186").Value = 311911
$ws.Range("C190").Value = 556
$ws.Range("D190").Value = 809382
$ws.Range("C192").Value = 85888
$ws.Range("D192").Value = 125786042
$ws.Range("C195").Value = 45307
$ws.Range("D195").Value = 65569267
$ws.Range("C198").Value = 4019
$ws.Range("D198").Value = 5654667
$ws.Range("C201").Value = 4320
$ws.Range("D201").Value = 5981557
$ws.Range("C204").Value = 94085
$ws.Range("D204").Value = 118431412
$ws.Range("C205").Value = 68
$ws.Range("D205").Value = 70104
$ws.Range("C209").Value = 508
$ws.Range("D209").Value = 742412
$ws.Range("C211").Value = 45479
$ws.Range("D211").Value = 66697294
$ws.Range("C213").Value = 10945
$ws.Range("D213").Value = 15750706
$ws.Range("C215").Value = 1687
$ws.Range("D215").Value = 2417347
$ws.Range("C216").Value = 8
$ws.Range("D216").Value = 10096
$ws.Range("C217").Value = 1948
$ws.Range("D217").Value = 2719580
$ws.Range("C218").Value = 225212
$ws.Range("D218").Value = 285059174
$ws.Range("C219").Value = 149
$ws.Range("D219").Value = 185356
$ws.Range("C220").Value = 225
$ws.Range("D220").Value = 324832
$ws.Range("C224").Value = 734
$ws.Range("D224").Value = 1080810
$ws.Range("C226").Value = 86433
$ws.Range("D226").Value = 126808751
$ws.Range("C229").Value = 57138
$ws.Range("D229").Value = 82914285
$ws.Range("C231").Value = 2082
$ws.Range("D231").Value = 2942645
$ws.Range("C234").Value = 3542
$ws.Range("D234").Value = 4972531
$ws.Range("C235").Value = 387829
$ws.Range("D235").Value = 509494954
$ws.Range("C236").Value = 70
$ws.Range("D236").Value = 90997
$ws.Range("C240").Value = 1152
$ws.Range("D240").Value = 1709551
$ws.Range("C242").Value = 179221
$ws.Range("D242").Value = 263875195
$ws.Range("C243").Value = 318
$ws.Range("D243").Value = 473790
$ws.Range("C245").Value = 149596
$ws.Range("D245").Value = 217607407
$ws.Range("C247").Value = 2239
$ws.Range("D247").Value = 3141648
$ws.Range("C249").Value = 4586
$ws.Range("D249").Value = 6471254
